$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Number of Party Members" (D1) goes from 1 to 3
$ws.Range("D1").Value = 3

# J12:J32 total-xp formulas now scale by the party size in $D$1
for ($r = 12; $r -le 32; $r++) {
    $ws.Range("J$r").Formula = "=C$r/D$r*`$D`$1"
}

# Move the selection/scroll position to match the saved view
$ws.Range("I18").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
